$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6612873673439026
$ws.Range("B1").Value = 1.15418529510498
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.081625461578369
$ws.Range("E1").Value = 1.090499639511108
